# "Add files via upload" - appends a second week ("Day 2") status block to the
# existing team work-progress sheet, re-using the same layout as the first
# block (rows 1-17) and leaves the selection positioned the way the author
# left it when they saved (cell E23, scrolled so row 7 is at the top).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Day 2" heading / date -------------------------------------------------
$ws.Range("C19").Value = "Day 2"
$ws.Range("C20").Value = "(14 Aug 2018)"

# --- Sushmitha S.H -----------------------------------------------------------
$ws.Range("A21").Value = "Sushmitha S.H"
$ws.Range("B22").Value = "Task Assigned"
$ws.Range("C22").Value = "Create POJO classes"
# Match the centered / wrap-text formatting used by the equivalent "Task
# Assigned" cells in the first block (C4 / C12) instead of re-creating a new
# style, by copying the format only.
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("B23").Value = "Status"
$ws.Range("C23").Value = "Completed"

# --- Sarath Nistala ----------------------------------------------------------
$ws.Range("A25").Value = "Sarath Nistala"
$ws.Range("B26").Value = "Task Assigned"
$ws.Range("C26").Value = "Module Implementation "
$ws.Range("B27").Value = "Status"
$ws.Range("C27").Value = "Started"

# --- Yamini Ambati -----------------------------------------------------------
$ws.Range("A29").Value = "Yamini Ambati"
$ws.Range("B30").Value = "Task Assigned"
$ws.Range("C30").Value = "Create POJO classes"
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C30").PasteSpecial(-4122) | Out-Null
$ws.Range("B31").Value = "Status"
$ws.Range("C31").Value = "Completed"

# --- Vikas Kumar Roy ---------------------------------------------------------
$ws.Range("A33").Value = "Vikas Kumar Roy"
$ws.Range("B34").Value = "Task Assigned"
$ws.Range("C34").Value = "Module Implementation"
$ws.Range("B35").Value = "Status"
$ws.Range("C35").Value = "Started"

# --- restore cursor / scroll position left behind in the saved file --------
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E23").Select() | Out-Null
